# Preface/Acknowledgements revisions per Frank's comments:
#  1. Re-number the TOC bookmark on the "Acknowledgements" heading
#     (_Toc5635009 -> _Toc5705284), keeping the same bookmark id/anchor.
#  2. Remove the closing "family and friends" acknowledgements paragraph.

$d = $word.ActiveDocument

# --- 1. Rename the second TOC bookmark on the heading -------------------
$oldBookmarkName = "_Toc5635009"
$newBookmarkName = "_Toc5705284"

if ($d.Bookmarks.Exists($oldBookmarkName)) {
    $bm = $d.Bookmarks($oldBookmarkName)
    $bmRange = $bm.Range
    $bm.Delete()
    $d.Bookmarks.Add($newBookmarkName, $bmRange) | Out-Null
}

# --- 2. Delete the trailing "family and friends" paragraph --------------
$target = "I" + [char]0x2019 + "d like to thank my family and friends for their good humor and support during this endeavor. I am especially grateful to my friends and colleagues, John DeVitis, David Masceri, Nick Romano and Dom Wirkijowski for their extensive help over these several years.  "

foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
        $para.Range.Delete()
        break
    }
}
